$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2 (anchor G2=5489, diff @@ -727)
$ws.Cells.Item(2, 8).Value = 301.33334
$ws.Cells.Item(2, 9).Value = 301.6
$ws.Cells.Item(2, 11).Value = 301.6
$ws.Cells.Item(2, 13).Value = -188.6
# row 5 (anchor G5=5503, diff @@ -880)
$ws.Cells.Item(5, 8).Value = 74.75
$ws.Cells.Item(5, 9).Value = 74.75
$ws.Cells.Item(5, 11).Value = 74.75
$ws.Cells.Item(5, 13).Value = 40.25
# row 9 (anchor G9=5487, diff @@ -1076)
$ws.Cells.Item(9, 8).Value = 170.1
$ws.Cells.Item(9, 9).Value = 158.71428
$ws.Cells.Item(9, 10).Value = 196.66667
$ws.Cells.Item(9, 11).Value = 158.71428
$ws.Cells.Item(9, 12).Value = 196.66667
$ws.Cells.Item(9, 13).Value = 10.28572
$ws.Cells.Item(9, 14).Value = -534.6666700000001
# row 18 (anchor G18=5471, diff @@ -1517)
$ws.Cells.Item(18, 8).Value = 1568.5
$ws.Cells.Item(18, 9).Value = 1568.5
$ws.Cells.Item(18, 11).Value = 1568.5
$ws.Cells.Item(18, 13).Value = -1284.5
# row 21 (anchor G21=2149, diff @@ -1667)
$ws.Cells.Item(21, 8).Value = 18000
$ws.Cells.Item(21, 9).Value = 18000
$ws.Cells.Item(21, 11).Value = 18000
$ws.Cells.Item(21, 13).Value = -17532
# row 23 (anchor G23=2149, diff @@ -1762)
$ws.Cells.Item(23, 8).Value = 18000
$ws.Cells.Item(23, 9).Value = 18000
$ws.Cells.Item(23, 11).Value = 18000
$ws.Cells.Item(23, 13).Value = -17766
# row 34 (anchor G34=2160, diff @@ -2298)
$ws.Cells.Item(34, 8).Value = 1977.4
$ws.Cells.Item(34, 9).Value = 1977.4
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1977.4
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -1774.4
$ws.Cells.Item(34, 14).Value = ""
# row 36 (anchor G36=2160, diff @@ -2399)
$ws.Cells.Item(36, 8).Value = 1977.4
$ws.Cells.Item(36, 9).Value = 1977.4
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 1977.4
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -1262.4
$ws.Cells.Item(36, 14).Value = ""
# row 55 (anchor G55=5517, diff @@ -3327)
$ws.Cells.Item(55, 8).Value = 146.55556
$ws.Cells.Item(55, 9).Value = 153
$ws.Cells.Item(55, 11).Value = 153
$ws.Cells.Item(55, 13).Value = 61
# row 62 (anchor G62=27781, diff @@ -3679)
$ws.Cells.Item(62, 8).Value = 4175.75
$ws.Cells.Item(62, 9).Value = 4175.75
$ws.Cells.Item(62, 11).Value = 4175.75
$ws.Cells.Item(62, 13).Value = -3551.75
# row 65 (anchor G65=27781, diff @@ -3823)
$ws.Cells.Item(65, 8).Value = 4175.75
$ws.Cells.Item(65, 9).Value = 4175.75
$ws.Cells.Item(65, 11).Value = 20878.75
$ws.Cells.Item(65, 13).Value = -17758.75
# row 76 (anchor G76=12602, diff @@ -4359)
$ws.Cells.Item(76, 8).Value = 5371.143
$ws.Cells.Item(76, 9).Value = 5371.143
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 5371.143
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -5056.143
$ws.Cells.Item(76, 14).Value = ""
# row 79 (anchor G79=12602, diff @@ -4503)
$ws.Cells.Item(79, 8).Value = 5371.143
$ws.Cells.Item(79, 9).Value = 5371.143
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 5371.143
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -4279.143
$ws.Cells.Item(79, 14).Value = ""
# row 86 (anchor G86=12603, diff @@ -4843)
$ws.Cells.Item(86, 8).Value = 3197.5
$ws.Cells.Item(86, 9).Value = 3398
$ws.Cells.Item(86, 11).Value = 3398
$ws.Cells.Item(86, 13).Value = -2275
# row 88 (anchor G88=12608, diff @@ -4944)
$ws.Cells.Item(88, 8).Value = 2419.8125
$ws.Cells.Item(88, 9).Value = 2608.75
$ws.Cells.Item(88, 10).Value = 2356.8333
$ws.Cells.Item(88, 11).Value = 2608.75
$ws.Cells.Item(88, 12).Value = 2356.8333
$ws.Cells.Item(88, 13).Value = -2202.75
$ws.Cells.Item(88, 14).Value = -3168.8333
# row 89 (anchor G89=12603, diff @@ -4996)
$ws.Cells.Item(89, 8).Value = 3197.5
$ws.Cells.Item(89, 9).Value = 3398
$ws.Cells.Item(89, 11).Value = 16990
$ws.Cells.Item(89, 13).Value = -11374
# row 91 (anchor G91=12608, diff @@ -5097)
$ws.Cells.Item(91, 8).Value = 2419.8125
$ws.Cells.Item(91, 9).Value = 2608.75
$ws.Cells.Item(91, 10).Value = 2356.8333
$ws.Cells.Item(91, 11).Value = 2608.75
$ws.Cells.Item(91, 12).Value = 2356.8333
$ws.Cells.Item(91, 13).Value = -1204.75
$ws.Cells.Item(91, 14).Value = -5164.8333
# row 96 (anchor G96=19894, diff @@ -5345)
$ws.Cells.Item(96, 8).Value = 1740
$ws.Cells.Item(96, 10).Value = 57
$ws.Cells.Item(96, 12).Value = 171
$ws.Cells.Item(96, 14).Value = -2917
# row 107 (anchor G107=27766, diff @@ -5896)
$ws.Cells.Item(107, 8).Value = 616.6667
$ws.Cells.Item(107, 9).Value = 625
$ws.Cells.Item(107, 10).Value = 600
$ws.Cells.Item(107, 11).Value = 625
$ws.Cells.Item(107, 12).Value = 600
$ws.Cells.Item(107, 13).Value = 1295
$ws.Cells.Item(107, 14).Value = -4440
# row 112 (anchor G112=27960, diff @@ -6138)
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).Value = ""
# row 125 (anchor G125=36228, diff @@ -6769)
$ws.Cells.Item(125, 8).Value = 1000.2222
$ws.Cells.Item(125, 9).Value = 1100
$ws.Cells.Item(125, 10).Value = 992.24
$ws.Cells.Item(125, 11).Value = 9900
$ws.Cells.Item(125, 12).Value = 8930.16
$ws.Cells.Item(125, 13).Value = -7440
$ws.Cells.Item(125, 14).Value = -13850.16
# row 127 (anchor G127=36114, diff @@ -6867)
$ws.Cells.Item(127, 8).Value = 2421.5
$ws.Cells.Item(127, 9).Value = 1305.8
$ws.Cells.Item(127, 11).Value = 3917.4
$ws.Cells.Item(127, 13).Value = 1042.6
# row 132 (anchor G132=44049, diff @@ -7112)
$ws.Cells.Item(132, 8).Value = 52636040
$ws.Cells.Item(132, 9).Value = 52636040
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 157908120
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -157905590
$ws.Cells.Item(132, 14).Value = ""
# row 137 (anchor G137=44013, diff @@ -7351)
$ws.Cells.Item(137, 8).Value = 2613
$ws.Cells.Item(137, 9).Value = 1419.6666
$ws.Cells.Item(137, 11).Value = 4258.9998
$ws.Cells.Item(137, 13).Value = -1708.9998
# row 138 (anchor G138=44169, diff @@ -7403)
$ws.Cells.Item(138, 8).Value = 2834.1904
$ws.Cells.Item(138, 10).Value = 3223.182
$ws.Cells.Item(138, 12).Value = 9669.545999999998
$ws.Cells.Item(138, 14).Value = -19949.546
# row 141 (anchor G141=44161, diff @@ -7553)
$ws.Cells.Item(141, 8).Value = 9676
$ws.Cells.Item(141, 9).Value = 9676
$ws.Cells.Item(141, 11).Value = 29028
$ws.Cells.Item(141, 13).Value = -23848

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2 (anchor G2=27713, diff @@ -7696)
$ws.Cells.Item(2, 8).Value = 1558.9546
$ws.Cells.Item(2, 9).Value = 1253.3334
$ws.Cells.Item(2, 11).Value = 1253.3334
$ws.Cells.Item(2, 13).Value = -1140.3334
# row 32 (anchor G32=44147, diff @@ -9160)
$ws.Cells.Item(32, 8).Value = 12624.12
$ws.Cells.Item(32, 9).Value = 12300.818
$ws.Cells.Item(32, 10).Value = 14995
$ws.Cells.Item(32, 11).Value = 12300.818
$ws.Cells.Item(32, 12).Value = 14995
$ws.Cells.Item(32, 13).Value = -12013.818
$ws.Cells.Item(32, 14).Value = -15569
# row 61 (anchor G61=43999, diff @@ -10560)
$ws.Cells.Item(61, 8).Value = 1566.1333
$ws.Cells.Item(61, 9).Value = 1587.4286
$ws.Cells.Item(61, 11).Value = 1587.4286
$ws.Cells.Item(61, 13).Value = -1375.4286
# row 62 (anchor G62=10719, diff @@ -10612)
$ws.Cells.Item(62, 8).Value = 99997
$ws.Cells.Item(62, 10).Value = 99997
$ws.Cells.Item(62, 12).Value = 99997
$ws.Cells.Item(62, 14).Value = -101245
# row 65 (anchor G65=10719, diff @@ -10759)
$ws.Cells.Item(65, 8).Value = 99997
$ws.Cells.Item(65, 10).Value = 99997
$ws.Cells.Item(65, 12).Value = 299991
$ws.Cells.Item(65, 14).Value = -306231
# row 74 (anchor G74=44000, diff @@ -11182)
$ws.Cells.Item(74, 8).Value = 2348.2
$ws.Cells.Item(74, 9).Value = 1849.4117
$ws.Cells.Item(74, 10).Value = 3408.125
$ws.Cells.Item(74, 11).Value = 1849.4117
$ws.Cells.Item(74, 12).Value = 3408.125
$ws.Cells.Item(74, 13).Value = -975.4117000000001
$ws.Cells.Item(74, 14).Value = -5156.125
# row 77 (anchor G77=44000, diff @@ -11329)
$ws.Cells.Item(77, 8).Value = 2348.2
$ws.Cells.Item(77, 9).Value = 1849.4117
$ws.Cells.Item(77, 10).Value = 3408.125
$ws.Cells.Item(77, 11).Value = 9247.058500000001
$ws.Cells.Item(77, 12).Value = 17040.625
$ws.Cells.Item(77, 13).Value = -4879.058500000001
$ws.Cells.Item(77, 14).Value = -25776.625
# row 97 (anchor G97=19941, diff @@ -12291)
$ws.Cells.Item(97, 8).Value = 2937.5
$ws.Cells.Item(97, 9).Value = 1370
$ws.Cells.Item(97, 10).Value = 6072.5
$ws.Cells.Item(97, 11).Value = 1370
$ws.Cells.Item(97, 12).Value = 6072.5
$ws.Cells.Item(97, 13).Value = -874
$ws.Cells.Item(97, 14).Value = -7064.5
# row 116 (anchor G116=27713, diff @@ -13189)
$ws.Cells.Item(116, 8).Value = 1558.9546
$ws.Cells.Item(116, 9).Value = 1253.3334
$ws.Cells.Item(116, 11).Value = 1253.3334
$ws.Cells.Item(116, 13).Value = 1040.6666
# row 122 (anchor G122=36168, diff @@ -13489)
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).Value = ""
# row 132 (anchor G132=43997, diff @@ -13961)
$ws.Cells.Item(132, 8).Value = 8994.6
$ws.Cells.Item(132, 9).Value = 9047.5
$ws.Cells.Item(132, 10).Value = 8945.77
$ws.Cells.Item(132, 11).Value = 27142.5
$ws.Cells.Item(132, 12).Value = 26837.31
$ws.Cells.Item(132, 13).Value = -24612.5
$ws.Cells.Item(132, 14).Value = -31897.31
# row 136 (anchor G136=43999, diff @@ -14157)
$ws.Cells.Item(136, 8).Value = 1566.1333
$ws.Cells.Item(136, 9).Value = 1587.4286
$ws.Cells.Item(136, 11).Value = 4762.2858
$ws.Cells.Item(136, 13).Value = -2212.2858

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3 (anchor G3=27713, diff @@ -14591)
$ws.Cells.Item(3, 8).Value = 1558.9546
$ws.Cells.Item(3, 9).Value = 1253.3334
$ws.Cells.Item(3, 11).Value = 1253.3334
$ws.Cells.Item(3, 13).Value = -1139.3334
# row 22 (anchor G22=5092, diff @@ -15522)
$ws.Cells.Item(22, 8).Value = 333
$ws.Cells.Item(22, 9).Value = 333
$ws.Cells.Item(22, 11).Value = 333
$ws.Cells.Item(22, 13).Value = -160
# row 62 (anchor G62=10586, diff @@ -17452)
$ws.Cells.Item(62, 8).Value = 65000
$ws.Cells.Item(62, 10).Value = 65000
$ws.Cells.Item(62, 12).Value = 65000
$ws.Cells.Item(62, 14).Value = -66372
# row 63 (anchor G63=10592, diff @@ -17498)
$ws.Cells.Item(63, 8).Value = 26135.5
$ws.Cells.Item(63, 9).Value = 2000
$ws.Cells.Item(63, 10).Value = 50271
$ws.Cells.Item(63, 11).Value = 2000
$ws.Cells.Item(63, 12).Value = 50271
$ws.Cells.Item(63, 13).Value = -1314
$ws.Cells.Item(63, 14).Value = -51643
# row 65 (anchor G65=10586, diff @@ -17596)
$ws.Cells.Item(65, 8).Value = 65000
$ws.Cells.Item(65, 10).Value = 65000
$ws.Cells.Item(65, 12).Value = 195000
$ws.Cells.Item(65, 14).Value = -201864
# row 66 (anchor G66=10592, diff @@ -17642)
$ws.Cells.Item(66, 8).Value = 26135.5
$ws.Cells.Item(66, 9).Value = 2000
$ws.Cells.Item(66, 10).Value = 50271
$ws.Cells.Item(66, 11).Value = 6000
$ws.Cells.Item(66, 12).Value = 150813
$ws.Cells.Item(66, 13).Value = -2568
$ws.Cells.Item(66, 14).Value = -157677
# row 105 (anchor G105=19947, diff @@ -19520)
$ws.Cells.Item(105, 8).Value = 2955.625
$ws.Cells.Item(105, 9).Value = 2892.1428
$ws.Cells.Item(105, 11).Value = 2892.1428
$ws.Cells.Item(105, 13).Value = -1145.1428
# row 107 (anchor G107=27706, diff @@ -19621)
$ws.Cells.Item(107, 8).Value = 3315.4443
$ws.Cells.Item(107, 9).Value = 3229.875
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 3229.875
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = -1309.875
$ws.Cells.Item(107, 14).Value = -7840
# row 134 (anchor G134=43998, diff @@ -20899)
$ws.Cells.Item(134, 8).Value = 1771.1818
$ws.Cells.Item(134, 9).Value = 1771.1818
$ws.Cells.Item(134, 11).Value = 5313.5454
$ws.Cells.Item(134, 13).Value = -2778.5454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22 (anchor G22=5367, diff @@ -22374)
$ws.Cells.Item(22, 8).Value = 1204.25
$ws.Cells.Item(22, 9).Value = 189.66667
$ws.Cells.Item(22, 11).Value = 189.66667
$ws.Cells.Item(22, 13).Value = 160.33333
# row 31 (anchor G31=44023, diff @@ -22821)
$ws.Cells.Item(31, 8).Value = 5163.593
$ws.Cells.Item(31, 9).Value = 7179.1665
$ws.Cells.Item(31, 11).Value = 7179.1665
$ws.Cells.Item(31, 13).Value = -6884.1665
# row 34 (anchor G34=44023, diff @@ -22974)
$ws.Cells.Item(34, 8).Value = 5163.593
$ws.Cells.Item(34, 9).Value = 7179.1665
$ws.Cells.Item(34, 11).Value = 7179.1665
$ws.Cells.Item(34, 13).Value = -6977.1665
# row 58 (anchor G58=44021, diff @@ -24156)
$ws.Cells.Item(58, 8).Value = 2445
$ws.Cells.Item(58, 9).Value = 2445
$ws.Cells.Item(58, 11).Value = 2445
$ws.Cells.Item(58, 13).Value = -2242
# row 64 (anchor G64=10610, diff @@ -24453)
$ws.Cells.Item(64, 8).Value = 95203.75
$ws.Cells.Item(64, 10).Value = 95203.75
$ws.Cells.Item(64, 12).Value = 95203.75
$ws.Cells.Item(64, 14).Value = -95699.75
# row 67 (anchor G67=10610, diff @@ -24594)
$ws.Cells.Item(67, 8).Value = 95203.75
$ws.Cells.Item(67, 10).Value = 95203.75
$ws.Cells.Item(67, 12).Value = 95203.75
$ws.Cells.Item(67, 14).Value = -96919.75
# row 86 (anchor G86=12584, diff @@ -25474)
$ws.Cells.Item(86, 8).Value = 19998
$ws.Cells.Item(86, 9).Value = 16997.6
$ws.Cells.Item(86, 11).Value = 16997.6
$ws.Cells.Item(86, 13).Value = -15874.6
# row 89 (anchor G89=12584, diff @@ -25621)
$ws.Cells.Item(89, 8).Value = 19998
$ws.Cells.Item(89, 9).Value = 16997.6
$ws.Cells.Item(89, 11).Value = 84988
$ws.Cells.Item(89, 13).Value = -79372
# row 105 (anchor G105=19928, diff @@ -26390)
$ws.Cells.Item(105, 8).Value = 1998.4
$ws.Cells.Item(105, 9).Value = 1950.0769
$ws.Cells.Item(105, 11).Value = 1950.0769
$ws.Cells.Item(105, 13).Value = -203.0769
# row 107 (anchor G107=27689, diff @@ -26491)
$ws.Cells.Item(107, 8).Value = 1010.2381
$ws.Cells.Item(107, 9).Value = 994.3570999999999
$ws.Cells.Item(107, 10).Value = 1042
$ws.Cells.Item(107, 11).Value = 994.3570999999999
$ws.Cells.Item(107, 12).Value = 1042
$ws.Cells.Item(107, 13).Value = 925.6429000000001
$ws.Cells.Item(107, 14).Value = -4882
# row 132 (anchor G132=44019, diff @@ -27695)
$ws.Cells.Item(132, 8).Value = 5734.25
$ws.Cells.Item(132, 10).Value = 5753.75
$ws.Cells.Item(132, 12).Value = 17261.25
$ws.Cells.Item(132, 14).Value = -22321.25
# row 134 (anchor G134=44020, diff @@ -27793)
$ws.Cells.Item(134, 8).Value = 1749.875
$ws.Cells.Item(134, 9).Value = 1749.875
$ws.Cells.Item(134, 11).Value = 5249.625
$ws.Cells.Item(134, 13).Value = -2714.625
# row 136 (anchor G136=44021, diff @@ -27891)
$ws.Cells.Item(136, 8).Value = 2445
$ws.Cells.Item(136, 9).Value = 2445
$ws.Cells.Item(136, 11).Value = 7335
$ws.Cells.Item(136, 13).Value = -4785

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2 (anchor G2=4847, diff @@ -28276)
$ws.Cells.Item(2, 8).Value = 818.8570999999999
$ws.Cells.Item(2, 9).Value = 39
$ws.Cells.Item(2, 10).Value = 1130.8
$ws.Cells.Item(2, 11).Value = 234
$ws.Cells.Item(2, 12).Value = 6784.799999999999
$ws.Cells.Item(2, 13).Value = -121
$ws.Cells.Item(2, 14).Value = -7010.799999999999
# row 4 (anchor G4=4650, diff @@ -28377)
$ws.Cells.Item(4, 8).Value = 1682900.2
$ws.Cells.Item(4, 9).Value = 17950.25
$ws.Cells.Item(4, 11).Value = 53850.75
$ws.Cells.Item(4, 13).Value = -53738.75
# row 23 (anchor G23=4858, diff @@ -29338)
$ws.Cells.Item(23, 8).Value = 868.8
$ws.Cells.Item(23, 9).Value = 995
$ws.Cells.Item(23, 10).Value = 837.25
$ws.Cells.Item(23, 11).Value = 2985
$ws.Cells.Item(23, 12).Value = 2511.75
$ws.Cells.Item(23, 13).Value = -2750
$ws.Cells.Item(23, 14).Value = -2981.75
# row 50 (anchor G50=4725, diff @@ -30667)
$ws.Cells.Item(50, 8).Value = 358.375
$ws.Cells.Item(50, 9).Value = 331
$ws.Cells.Item(50, 11).Value = 993
$ws.Cells.Item(50, 13).Value = -512
# row 53 (anchor G53=4725, diff @@ -30820)
$ws.Cells.Item(53, 8).Value = 358.375
$ws.Cells.Item(53, 9).Value = 331
$ws.Cells.Item(53, 11).Value = 993
$ws.Cells.Item(53, 13).Value = -512
# row 113 (anchor G113=27843, diff @@ -33742)
$ws.Cells.Item(113, 8).Value = 1077.8541
$ws.Cells.Item(113, 9).Value = 1077.8541
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 3233.5623
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -1063.5623
$ws.Cells.Item(113, 14).Value = ""
# row 119 (anchor G119=27873, diff @@ -34039)
$ws.Cells.Item(119, 8).Value = 5001166
$ws.Cells.Item(119, 9).Value = 5001166
$ws.Cells.Item(119, 11).Value = 15003498
$ws.Cells.Item(119, 13).Value = -14998660
# row 131 (anchor G131=36060, diff @@ -34636)
$ws.Cells.Item(131, 8).Value = 1471.4828
$ws.Cells.Item(131, 10).Value = 1427.1228
$ws.Cells.Item(131, 12).Value = 4281.3684
$ws.Cells.Item(131, 14).Value = -14361.3684
# row 138 (anchor G138=44105, diff @@ -34991)
$ws.Cells.Item(138, 8).Value = 8451.388999999999
$ws.Cells.Item(138, 9).Value = 8451.388999999999
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 25354.167
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = -20214.167
$ws.Cells.Item(138, 14).Value = ""
# row 139 (anchor G139=44102, diff @@ -35043)
$ws.Cells.Item(139, 8).Value = 7361.16
$ws.Cells.Item(139, 9).Value = 4805.8
$ws.Cells.Item(139, 11).Value = 14417.4
$ws.Cells.Item(139, 13).Value = -9277.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 2 (anchor G2=5062, diff @@ -35287)
$ws.Cells.Item(2, 8).Value = 55.625
$ws.Cells.Item(2, 9).Value = 3.3333333
$ws.Cells.Item(2, 10).Value = 87
$ws.Cells.Item(2, 11).Value = 3.3333333
$ws.Cells.Item(2, 12).Value = 87
$ws.Cells.Item(2, 13).Value = 109.6666667
$ws.Cells.Item(2, 14).Value = -313
# row 95 (anchor G95=18235, diff @@ -39751)
$ws.Cells.Item(95, 8).Value = 25257.5
$ws.Cells.Item(95, 10).Value = 25257.5
$ws.Cells.Item(95, 12).Value = 25257.5
$ws.Cells.Item(95, 14).Value = -30749.5
# row 97 (anchor G97=19940, diff @@ -39852)
$ws.Cells.Item(97, 8).Value = 876.9091
$ws.Cells.Item(97, 9).Value = 694.6
$ws.Cells.Item(97, 11).Value = 694.6
$ws.Cells.Item(97, 13).Value = -198.6
# row 102 (anchor G102=36169, diff @@ -40106)
$ws.Cells.Item(102, 8).Value = 2503.4285
$ws.Cells.Item(102, 9).Value = 1815.1111
$ws.Cells.Item(102, 10).Value = 6633.3335
$ws.Cells.Item(102, 11).Value = 1815.1111
$ws.Cells.Item(102, 12).Value = 6633.3335
$ws.Cells.Item(102, 13).Value = -193.1111000000001
$ws.Cells.Item(102, 14).Value = -9877.333500000001
# row 104 (anchor G104=18666, diff @@ -40207)
$ws.Cells.Item(104, 8).Value = 83490.86
$ws.Cells.Item(104, 10).Value = 83490.86
$ws.Cells.Item(104, 12).Value = 83490.86
$ws.Cells.Item(104, 14).Value = -90478.86
# row 107 (anchor G107=27802, diff @@ -40351)
$ws.Cells.Item(107, 8).Value = 2831.111
$ws.Cells.Item(107, 10).Value = 4396.6
$ws.Cells.Item(107, 12).Value = 4396.6
$ws.Cells.Item(107, 14).Value = -8236.6
# row 122 (anchor G122=36182, diff @@ -41059)
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = ""
$ws.Cells.Item(122, 14).Value = ""
# row 126 (anchor G126=36184, diff @@ -41252)
$ws.Cells.Item(126, 8).Value = 4499
$ws.Cells.Item(126, 9).Value = 4499
$ws.Cells.Item(126, 11).Value = 13497
$ws.Cells.Item(126, 13).Value = -11027
# row 132 (anchor G132=44008, diff @@ -41534)
$ws.Cells.Item(132, 8).Value = 3956.7693
$ws.Cells.Item(132, 9).Value = 2744.6667
$ws.Cells.Item(132, 10).Value = 4995.7144
$ws.Cells.Item(132, 11).Value = 8234.000100000001
$ws.Cells.Item(132, 12).Value = 14987.1432
$ws.Cells.Item(132, 13).Value = -5704.000100000001
$ws.Cells.Item(132, 14).Value = -20047.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 26 (anchor G26=3559, diff @@ -43282)
$ws.Cells.Item(26, 8).Value = 10000
$ws.Cells.Item(26, 10).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 14).Value = -10590
# row 40 (anchor G40=36248, diff @@ -43968)
$ws.Cells.Item(40, 8).Value = 125000776
$ws.Cells.Item(40, 9).Value = 142857940
$ws.Cells.Item(40, 10).Value = 700
$ws.Cells.Item(40, 11).Value = 142857940
$ws.Cells.Item(40, 12).Value = 700
$ws.Cells.Item(40, 13).Value = -142857804
$ws.Cells.Item(40, 14).Value = -972
# row 56 (anchor G56=3668, diff @@ -44728)
$ws.Cells.Item(56, 8).Value = 9849.75
$ws.Cells.Item(56, 9).Value = 9849.75
$ws.Cells.Item(56, 11).Value = 9849.75
$ws.Cells.Item(56, 13).Value = -9158.75
# row 62 (anchor G62=10740, diff @@ -45019)
$ws.Cells.Item(62, 8).Value = 20100
$ws.Cells.Item(62, 9).Value = 16000
$ws.Cells.Item(62, 10).Value = 24200
$ws.Cells.Item(62, 11).Value = 16000
$ws.Cells.Item(62, 12).Value = 24200
$ws.Cells.Item(62, 13).Value = -15376
$ws.Cells.Item(62, 14).Value = -25448
# row 65 (anchor G65=10740, diff @@ -45166)
$ws.Cells.Item(65, 8).Value = 20100
$ws.Cells.Item(65, 9).Value = 16000
$ws.Cells.Item(65, 10).Value = 24200
$ws.Cells.Item(65, 11).Value = 48000
$ws.Cells.Item(65, 12).Value = 72600
$ws.Cells.Item(65, 13).Value = -44880
$ws.Cells.Item(65, 14).Value = -78840
# row 76 (anchor G76=10742, diff @@ -45705)
$ws.Cells.Item(76, 8).Value = 16596.666
$ws.Cells.Item(76, 10).Value = 16596.666
$ws.Cells.Item(76, 12).Value = 16596.666
$ws.Cells.Item(76, 14).Value = -17272.666
# row 79 (anchor G79=10742, diff @@ -45852)
$ws.Cells.Item(79, 8).Value = 16596.666
$ws.Cells.Item(79, 10).Value = 16596.666
$ws.Cells.Item(79, 12).Value = 16596.666
$ws.Cells.Item(79, 14).Value = -18936.666
# row 100 (anchor G100=19995, diff @@ -46851)
$ws.Cells.Item(100, 8).Value = 2296
$ws.Cells.Item(100, 9).Value = 1995.3334
$ws.Cells.Item(100, 11).Value = 1995.3334
$ws.Cells.Item(100, 13).Value = -1454.3334
# row 122 (anchor G122=36247, diff @@ -47902)
$ws.Cells.Item(122, 8).Value = 5599
$ws.Cells.Item(122, 9).Value = 4998.75
$ws.Cells.Item(122, 11).Value = 14996.25
$ws.Cells.Item(122, 13).Value = -12546.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62 (anchor G62=12589, diff @@ -51856)
$ws.Cells.Item(62, 8).Value = 3994
$ws.Cells.Item(62, 9).Value = 3994
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 3994
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -3370
$ws.Cells.Item(62, 14).Value = ""
# row 65 (anchor G65=12589, diff @@ -52003)
$ws.Cells.Item(65, 8).Value = 3994
$ws.Cells.Item(65, 9).Value = 3994
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 19970
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -16850
$ws.Cells.Item(65, 14).Value = ""
# row 105 (anchor G105=18710, diff @@ -53939)
$ws.Cells.Item(105, 8).Value = 13990.223
$ws.Cells.Item(105, 9).Value = 19800
$ws.Cells.Item(105, 10).Value = 13264
$ws.Cells.Item(105, 11).Value = 19800
$ws.Cells.Item(105, 12).Value = 13264
$ws.Cells.Item(105, 13).Value = -16306
$ws.Cells.Item(105, 14).Value = -20252
# row 122 (anchor G122=36208, diff @@ -54757)
$ws.Cells.Item(122, 8).Value = 1410
$ws.Cells.Item(122, 9).Value = 1410
$ws.Cells.Item(122, 11).Value = 4230
$ws.Cells.Item(122, 13).Value = -1780
# row 132 (anchor G132=44029, diff @@ -55235)
$ws.Cells.Item(132, 8).Value = 26670.334
$ws.Cells.Item(132, 9).Value = 26670.334
$ws.Cells.Item(132, 11).Value = 80011.00199999999
$ws.Cells.Item(132, 13).Value = -77481.00199999999
# row 136 (anchor G136=44031, diff @@ -55434)
$ws.Cells.Item(136, 8).Value = 51059.133
$ws.Cells.Item(136, 9).Value = 50332.332
$ws.Cells.Item(136, 11).Value = 150996.996
$ws.Cells.Item(136, 13).Value = -148446.996

